$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the second sheet: "LandTaxAssessments" -> "Assessments"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Assessments"

# ---------------------------------------------------------------------------
# 2. CustomerData sheet: split ClientName into First Name / Last Name columns
# ---------------------------------------------------------------------------
# Insert two new (blank) columns right after column B (ClientName); they
# inherit column B's cell formatting automatically.
$ws1.Columns.Item(3).EntireColumn.Insert()
$ws1.Columns.Item(3).EntireColumn.Insert()

# Match column B's on-screen width for the two new columns.
$ws1.Range("C1:D1").ColumnWidth = $ws1.Columns.Item(2).ColumnWidth

# Headers
$ws1.Cells.Item(1, 3).Value = "First Name"
$ws1.Cells.Item(1, 4).Value = "Last Name"

# Existing rows: split the full name already in column B
$ws1.Cells.Item(2, 3).Value = "Sarah"
$ws1.Cells.Item(2, 4).Value = "Jones"

$ws1.Cells.Item(3, 3).Value = "Susan"
$ws1.Cells.Item(3, 4).Value = "Williams"

$ws1.Cells.Item(4, 3).Value = "Jane"
$ws1.Cells.Item(4, 4).Value = "Smith"

# ---------------------------------------------------------------------------
# 3. CustomerData sheet: add new customer row (John Citizen)
# ---------------------------------------------------------------------------
# A5, E5 and F5 carry the same plain style as the rows above (s="1"); reuse
# that existing style via a format-only paste instead of touching Font/
# NumberFormat directly (which would otherwise mint brand-new style entries).
$ws1.Range("A4").Copy()
$ws1.Range("A5").PasteSpecial(-4122)
$ws1.Range("E4:F4").Copy()
$ws1.Range("E5:F5").PasteSpecial(-4122)

$ws1.Cells.Item(5, 1).Value = 52457
$ws1.Cells.Item(5, 2).Value = "John Citizen"
$ws1.Cells.Item(5, 3).Value = "John"
$ws1.Cells.Item(5, 4).Value = "Citizen"
$ws1.Cells.Item(5, 5).Value = "CENTENNIAL PLAZA 260 ELIZABETH ST SURRY HILLS 2010 NSW"
$ws1.Cells.Item(5, 6).Value = 619876543223
$ws1.Cells.Item(5, 7).Value = "johncitizen123@gmail.com"

$ws1.Range("H4").Copy()
$ws1.Range("H5").PasteSpecial(-4122)
$ws1.Cells.Item(5, 8).Value = 27450

# ---------------------------------------------------------------------------
# 4. Assessments sheet: rename CorrespondenceID header + highlight it
# ---------------------------------------------------------------------------
$ws2.Cells.Item(1, 2).Value = "CorrespondenceID (AssessmentID from CM)"
$ws2.Cells.Item(1, 2).Font.Bold = $true
$ws2.Cells.Item(1, 2).Font.Color = 255

# ---------------------------------------------------------------------------
# 5. Selection: leave CustomerData's cursor at A6, restore Assessments as the
#    active tab with its cursor at A2.
# ---------------------------------------------------------------------------
$ws1.Range("A6").Select()
$ws2.Range("A2").Select()
